$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

$row = 61

$ws.Cells.Item($row, 1).Value = "Dit is echt niet netjes verlopen zo"
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #20: Dit is echt niet netjes verlopen zo"
$ws.Cells.Item($row, 4).Value = "Overig"
$ws.Cells.Item($row, 5).Value = "Beste klant,`nBedankt voor uw bericht. Kunt u wat meer details geven over wat er precies niet netjes is verlopen? Met deze informatie kunnen we het probleem beter begrijpen en u een passende oplossing bieden.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$ws.Cells.Item($row, 6).Value = "2025-08-05 20:05:12"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Nee"
$ws.Cells.Item($row, 9).Value = "Ja"
$ws.Cells.Item($row, 10).Value = "Nee"

# The engine auto-calculates a custom row height when multi-line text is
# assigned via .Value; re-running AutoFit resets it back to the default
# (un-customized) row height, matching the target workbook.
$ws.Rows($row).EntireRow.AutoFit()

# Extend the conditional-formatting ranges so they keep covering the full
# data range (2:60 -> 2:61) for every formatted column.
$ws.Range("D2:D60").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D61"))
$ws.Range("G2:G60").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G61"))
$ws.Range("H2:H60").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H61"))
$ws.Range("I2:I60").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I61"))
$ws.Range("J2:J60").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J61"))

# Dashboard summary count for the "Overig" category increased by one.
$dash.Cells.Item(3, 2).Value = 14
